# Auto-generated edit script applying cryptos.xlsx price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.383.12"
$ws.Range("E2").Value = "  -2.04%  "
$ws.Range("D3").Value = "3.378.25"
$ws.Range("E3").Value = "  -2.76%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'594.64"
$ws.Range("E5").Value = "  -1.58%  "
$ws.Range("D6").Value = "'141.27"
$ws.Range("E6").Value = "  -4.85%  "
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("D8").Value = "3.377.71"
$ws.Range("E8").Value = "  -2.71%  "
$ws.Range("D9").Value = "'0.467"
$ws.Range("E9").Value = "  -3.16%  "
$ws.Range("D10").Value = "'7.90"
$ws.Range("E10").Value = "  +4.68%  "
$ws.Range("D11").Value = "'0.132"
$ws.Range("E11").Value = "  -6.99%  "
$ws.Range("D12").Value = "'0.403"
$ws.Range("E12").Value = "  -5.06%  "
$ws.Range("D13").Value = "3.955.61"
$ws.Range("E13").Value = "  -2.63%  "
$ws.Range("D14").Value = "'0.0000199"
$ws.Range("E14").Value = "  -7.43%  "
$ws.Range("D15").Value = "'29.44"
$ws.Range("E15").Value = "  -7.24%  "
$ws.Range("D17").Value = "65.454.43"
$ws.Range("E17").Value = "  -2.08%  "
$ws.Range("D18").Value = "3.383.93"
$ws.Range("E18").Value = "  -2.57%  "
$ws.Range("D19").Value = "'10.27"
$ws.Range("E19").Value = "  +1.39%  "
$ws.Range("D20").Value = "'6.07"
$ws.Range("E20").Value = "  -6.13%  "
$ws.Range("D21").Value = "'14.51"
$ws.Range("E21").Value = "  -6.13%  "
$ws.Range("D22").Value = "'411.67"
$ws.Range("E22").Value = "  -6.40%  "
$ws.Range("D23").Value = "'0.576"
$ws.Range("E23").Value = "  -5.98%  "
$ws.Range("D24").Value = "'76.86"
$ws.Range("E24").Value = "  -3.09%  "
$ws.Range("E25").Value = "  -0.02%  "
$ws.Range("D26").Value = "3.520.97"
$ws.Range("E26").Value = "  -2.53%  "
$ws.Range("E27").Value = "  -10.14%  "
$ws.Range("D28").Value = "'9.16"
$ws.Range("E28").Value = "  -6.36%  "
$ws.Range("E29").Value = "  -7.92%  "
$ws.Range("D30").Value = "'2.40"
$ws.Range("E30").Value = "  -3.35%  "
$ws.Range("D31").Value = "'0.999"
$ws.Range("E31").Value = "  -0.19%  "
$ws.Range("D32").Value = "'0.159"
$ws.Range("E32").Value = "  -5.20%  "
$ws.Range("E33").Value = "  -8.88%  "
$ws.Range("D34").Value = "'24.23"
$ws.Range("E34").Value = "  -4.75%  "
$ws.Range("D35").Value = "3.379.18"
$ws.Range("E35").Value = "  -2.46%  "
$ws.Range("E36").Value = "  -0.07%  "
$ws.Range("B37").Value = "NEARProtocol"
$ws.Range("C37").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D37").Value = "'5.50"
$ws.Range("E37").Value = "  -9.33%  "
$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").Value = "'1.67"
$ws.Range("E38").Value = "  -7.72%  "
$ws.Range("D39").Value = "'1.00"
$ws.Range("E39").Value = "  +0.17%  "
$ws.Range("D40").Value = "'7.45"
$ws.Range("E40").Value = "  -6.18%  "
$ws.Range("D41").Value = "'168.32"
$ws.Range("E41").Value = "  -4.56%  "
$ws.Range("E42").Value = "  -4.63%  "
$ws.Range("D43").Value = "'0.864"
$ws.Range("E43").Value = "  -2.71%  "
$ws.Range("E44").Value = "  -7.92%  "
$ws.Range("E45").Value = "  -11.05%  "
$ws.Range("D46").Value = "'45.32"
$ws.Range("E46").Value = "  -2.01%  "
$ws.Range("D47").Value = "'26.34"
$ws.Range("E47").Value = "  -9.74%  "
$ws.Range("E48").Value = "  -5.81%  "
$ws.Range("D50").Value = "'2.24"
$ws.Range("E50").Value = "  -8.98%  "
$ws.Range("D51").Value = "'0.911"
$ws.Range("E51").Value = "  -7.60%  "
